# Update to Salinas Design - apply review changes from Infiniti Solutions.
#
# Work from the bottom of the document upward so that paragraph indices
# for content we haven't touched yet are never disturbed by an earlier
# paragraph deletion.

$d = $word.ActiveDocument

function Replace-InParagraph([int]$index, [string]$old, [string]$new) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# Remove the hidden "_GoBack" bookmark that used to sit on the Version line
# -- it gets relocated below to the end of the page-7 bullet. Do this first
# so the name is free before we re-add it.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# Page 7 bullet: append "(layout update)" and plant the relocated "_GoBack"
# bookmark right after the new text, before the paragraph mark (mirrors the
# bookmark that used to sit on the "Version:" line). A trailing sentinel
# character gives the bookmark a real (non-paragraph-end) anchor; once the
# bookmark is placed around it the sentinel is deleted, leaving the
# bookmark collapsed right after the visible text.
# ---------------------------------------------------------------------------
Replace-InParagraph 29 "Add P20, a 1x2 header; connection as shown in schematic" "Add P20, a 1x2 header; connection as shown in schematic (layout update)Z"

$tail = $d.Paragraphs(29).Range.Duplicate
$tail.MoveEnd(1, -1)
$sentinel = $d.Range($tail.End - 1, $tail.End)
$d.Bookmarks.Add("_GoBack", $sentinel) | Out-Null
$sentinel.Text = ""

# ---------------------------------------------------------------------------
# Page 6 bullets: the J8 placement note becomes the microUSB replacement
# bullet; the old "optional" / "need to discuss" bullets that used to follow
# it are dropped (the U18 DNP bullet stays put).
# ---------------------------------------------------------------------------
$d.Paragraphs(26).Range.Delete()
$d.Paragraphs(25).Range.Delete()
Replace-InParagraph 24 "Rev 1.0 placement for J8 is not correct. The connector is placed too far in; there is no clearance for connecting cable. " "Replace J8 with microUSB connector (layout update)"

# ---------------------------------------------------------------------------
# Page 5 bullet: reword the reason for the ground-connection fix.
# ---------------------------------------------------------------------------
Replace-InParagraph 22 "Add ground connection for J6 pin 2 (error from previous design)" "Add ground connection for J6 pin 2 (layout update)"

# ---------------------------------------------------------------------------
# Page 4 bullets: "Remove ..." becomes "Change ... to 0 Ohm ..."; the
# VCCB/VCCA connection bullets are repurposed to carry the R35/R36/R37 and
# TXB0104DR BOM-change text, and the bullets that used to carry that text
# (P5 pin 1, Replace R35.., Replace TXB..) are removed.
# ---------------------------------------------------------------------------
$d.Paragraphs(20).Range.Delete()
$d.Paragraphs(19).Range.Delete()
$d.Paragraphs(18).Range.Delete()
Replace-InParagraph 17 "Connect VCCA to VDDIO" "Change TXB0104DR to TXS0104EDR (BOM change)"
Replace-InParagraph 16 "Connect VCCB to 3P3V" "Change R35, R36 and R37 from 10KOhm to 1KOhm (BOM change)"
Replace-InParagraph 15 "Remove R30, R31 and R34 (BOM change)" "Change  R30, R31 and R34 to 0 Ohm (BOM change)"

# ---------------------------------------------------------------------------
# Page 3 bullets: same "Remove -> Change ... to 0 Ohm" treatment; the two
# VCCB/VCCA connection bullets that followed are dropped so the BOM text
# flows straight into the TXB0104DR replacement bullet.
# ---------------------------------------------------------------------------
$d.Paragraphs(12).Range.Delete()
$d.Paragraphs(11).Range.Delete()
Replace-InParagraph 10 "Remove R23, R24, R25, R26, R27 and R28 (BOM change)" "Change R23, R24, R25, R26, R27 and R28 to 0 Ohm (BOM change)"

# ---------------------------------------------------------------------------
# Header block: bump the date and version number.
# ---------------------------------------------------------------------------
Replace-InParagraph 2 "Version: 1.0" "Version: 1.1"
Replace-InParagraph 1 "Last Update: September 2, 2020" "Last Update: September 4, 2020"
